# Atualizacao sabao em barra
# Inserts a new product row ("SABAO EM BARRA - 200G - 200G" / S010046)
# into Planilha1 right before "SABONETE LIQUIDO - 5000ML" (old row 82),
# shifting the remaining rows down by one, then fixes up the dependent
# ranges (defined name FilterDatabase, conditional formatting range) and
# the saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 82 (copies formatting/styles from the row
#    above, matching the diff's s="3"/s="4" styling for columns A-D).
$ws.Rows("82:82").Insert()

# 2. Populate the new row with the new product data.
$ws.Range("A82").Value2 = "SABAO EM BARRA - 200G - 200G"
$ws.Range("B82").Value2 = "UN"
$ws.Range("C82").Value2 = "S010046"
$ws.Range("D82").Value2 = 51

# 3. Extend the _FilterDatabase defined name to cover the new last row.
$fdb = $wb.Names.Item("Planilha1!_FilterDatabase")
$fdb.RefersTo = "=Planilha1!`$A`$1:`$D`$97"

# 4. Extend the conditional formatting range (A2:D96 -> A2:D97) that
#    highlights non-blank rows.
$cfRange = $ws.Range("A2:D97")
$cf = $ws.Cells.FormatConditions.Item(2)
$cf.ModifyAppliesToRange($cfRange)

# 5. Restore the saved selection/scroll position shown in the diff.
$ws.Range("G81").Select()
